$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 45

# Columns A-L are plain text in the source data (even the numeric-looking
# ones), so force text entry with a leading apostrophe the way Excel's UI
# does, then restore the default "Normal" style so no stray number format
# sticks to the cell (matches the rest of the sheet, which carries no
# explicit style on data rows).
$ws.Cells.Item($row, 1).Value = "'-468"
$ws.Cells.Item($row, 2).Value = "'6/4/2025"
$ws.Cells.Item($row, 3).Value = "Chile 2305"
$ws.Cells.Item($row, 4).Value = "'3"
$ws.Cells.Item($row, 5).Value = "'807168187"
$ws.Cells.Item($row, 6).Value = "GESTION TELECENTRO"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Columna con base podrida nodo propio  telecentro  con rienda pique   sin riesgo de caida al 0 406 25 "
$ws.Cells.Item($row, 9).Value = "'1"
$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Nodo TLC"
$ws.Cells.Item($row, 12).Value = "Terminal"
$ws.Cells.Item($row, 13).Value = -58.399193
$ws.Cells.Item($row, 14).Value = -34.617418

$ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 12)).Style = "Normal"
